$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 210, shifting existing rows 210:222 down to 211:223
$ws.Rows.Item(210).Insert()

$ws.Cells.Item(210, 1).Value = 9
$ws.Cells.Item(210, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(210, 3).Value = "Metropolitana"
$ws.Cells.Item(210, 4).Value = 44610
$ws.Cells.Item(210, 5).Value = 13
$ws.Cells.Item(210, 6).Value = 300000001
$ws.Cells.Item(210, 7).Value = "Rabanito"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 5200
$ws.Cells.Item(210, 11).Value = 2500
$ws.Cells.Item(210, 12).Value = 3000
$ws.Cells.Item(210, 13).Value = 2750
$ws.Cells.Item(210, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(210, 15).Value = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value = 28
$ws.Cells.Item(210, 17).Value = 100
$ws.Cells.Item(210, 18).Value = "Hortaliza"
